# Weekly update: insert 3 new "Pimiento" price rows (Feria Lagunitas de
# Puerto Montt, week of 2023-01-13 / serial 44939) at the top of the data
# block that starts at row 936, pushing the existing historical rows down
# by three (the sheet had rows 2..986 of data; it now has rows 2..989).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three blank rows before the current row 936; everything that was
# at 936.. shifts down to 939.. (and the last three existing rows end up
# at 987-989), matching the diff exactly.
$ws.Rows("936:938").Insert()

# Common/static columns shared by every data row in this sheet.
$mercadoId = 4
$mercado = "Feria Lagunitas de Puerto Montt"
$region = "Los Lagos"
$codreg = 10
$categoriaId = 100112002
$categoria = "Pimiento"
$clasificacion = "Hortaliza"
$fecha = 44939

# Row 936: Zafiro rojo / Primera
$ws.Cells.Item(936, 1).Value = $mercadoId
$ws.Cells.Item(936, 2).Value = $mercado
$ws.Cells.Item(936, 3).Value = $region
$ws.Cells.Item(936, 4).Value = $fecha
$ws.Cells.Item(936, 5).Value = $codreg
$ws.Cells.Item(936, 6).Value = $categoriaId
$ws.Cells.Item(936, 7).Value = $categoria
$ws.Cells.Item(936, 8).Value = "Zafiro rojo"
$ws.Cells.Item(936, 9).Value = "Primera"
$ws.Cells.Item(936, 10).Value = 150
$ws.Cells.Item(936, 11).Value = 18000
$ws.Cells.Item(936, 12).Value = 18000
$ws.Cells.Item(936, 13).Value = 18000
$ws.Cells.Item(936, 14).Value = "$/caja 15 kilos"
$ws.Cells.Item(936, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(936, 16).Value = 1200
$ws.Cells.Item(936, 17).Value = 15
$ws.Cells.Item(936, 18).Value = $clasificacion

# Row 937: Zafiro rojo / Primera
$ws.Cells.Item(937, 1).Value = $mercadoId
$ws.Cells.Item(937, 2).Value = $mercado
$ws.Cells.Item(937, 3).Value = $region
$ws.Cells.Item(937, 4).Value = $fecha
$ws.Cells.Item(937, 5).Value = $codreg
$ws.Cells.Item(937, 6).Value = $categoriaId
$ws.Cells.Item(937, 7).Value = $categoria
$ws.Cells.Item(937, 8).Value = "Zafiro rojo"
$ws.Cells.Item(937, 9).Value = "Primera"
$ws.Cells.Item(937, 10).Value = 120
$ws.Cells.Item(937, 11).Value = 25000
$ws.Cells.Item(937, 12).Value = 25000
$ws.Cells.Item(937, 13).Value = 25000
$ws.Cells.Item(937, 14).Value = "$/caja 18 kilos"
$ws.Cells.Item(937, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(937, 16).Value = 1389
$ws.Cells.Item(937, 17).Value = 18
$ws.Cells.Item(937, 18).Value = $clasificacion

# Row 938: Zafiro verde / Primera
$ws.Cells.Item(938, 1).Value = $mercadoId
$ws.Cells.Item(938, 2).Value = $mercado
$ws.Cells.Item(938, 3).Value = $region
$ws.Cells.Item(938, 4).Value = $fecha
$ws.Cells.Item(938, 5).Value = $codreg
$ws.Cells.Item(938, 6).Value = $categoriaId
$ws.Cells.Item(938, 7).Value = $categoria
$ws.Cells.Item(938, 8).Value = "Zafiro verde"
$ws.Cells.Item(938, 9).Value = "Primera"
$ws.Cells.Item(938, 10).Value = 250
$ws.Cells.Item(938, 11).Value = 20000
$ws.Cells.Item(938, 12).Value = 20000
$ws.Cells.Item(938, 13).Value = 20000
$ws.Cells.Item(938, 14).Value = "$/caja 15 kilos"
$ws.Cells.Item(938, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(938, 16).Value = 1333
$ws.Cells.Item(938, 17).Value = 15
$ws.Cells.Item(938, 18).Value = $clasificacion
